$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New quarterly row (row 60) appended after the last data row (59).
# Column A holds a date-like label that must stay plain text (matching the
# existing "01-MM-YYYY" labels in the rest of the column), so we briefly mark
# the cell as text before typing it in, then restore the default ("Normal")
# cell style so no stray formatting is left behind on the cell.
$ws.Range("A60").NumberFormat = "@"
$ws.Range("A60").Value = "01-07-2021"
$ws.Range("A60").Style = "Normal"

$ws.Range("B60").Value = 7331
$ws.Range("D60").Value = 3633
$ws.Range("E60").Value = -36
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 1427
$ws.Range("H60").Value = 2308
